$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCD")

# Insert a new blank column before D; existing D:K data shifts to E:L
$ws.Columns("D").Insert()

# Copy number formats from column E into the new column D for the data blocks
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate new column D with the latest (most recent) period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 21025200
$ws.Range("D9").Value = 10239200
$ws.Range("D10").Value = 10786000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -72400
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 12202600
$ws.Range("D18").Value = 8822600
$ws.Range("D20").Value = -25300
$ws.Range("D21").Value = 10279300
$ws.Range("D22").Value = 981200
$ws.Range("D23").Value = 7816100
$ws.Range("D24").Value = 1816800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 5999300
$ws.Range("D27").Value = 5999300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -75000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 25300
$ws.Range("D33").Value = 5924300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 5924300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 866000
$ws.Range("D42").Value = 171600
$ws.Range("D43").Value = 2441500
$ws.Range("D44").Value = 51100
$ws.Range("D45").Value = 523000
$ws.Range("D46").Value = 4053200
$ws.Range("D47").Value = 1202800
$ws.Range("D48").Value = 22842700
$ws.Range("D49").Value = 2331500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2381000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 32811200
$ws.Range("D57").Value = 1207900
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 1765600
$ws.Range("D60").Value = 2973500
$ws.Range("D61").Value = 31075300
$ws.Range("D62").Value = 5020800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 39069600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 50487000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -6258400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 5924300
$ws.Range("D83").Value = 1482000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 6966700
$ws.Range("D91").Value = -2741700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2455100
$ws.Range("D96").Value = -3255900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -5949600
$ws.Range("D101").Value = -159800
$ws.Range("D102").Value = -1597800
